# Moore-Heather timesheet: submit Week 11 timesheet entries
# - Adds two new entries (rows 6-7) to "Week 10"
# - Adds four new entries (rows 2-5) to "Week 11"
# - Moves the active/selected tab from "Week 10" to "Week 11"

$wb = $excel.ActiveWorkbook

$dateFmt = "m/d/yy"
$timeFmt = "h:mm AM/PM"

# ---------------------------------------------------------------------
# Week 10 - add rows 6 and 7
# ---------------------------------------------------------------------
$ws10 = $wb.Worksheets.Item("Week 10")

# Row 6
$ws10.Range("A6").Value = 41712
$ws10.Range("A6").NumberFormat = $dateFmt
$ws10.Range("B6").Value = 0.60416666666666663
$ws10.Range("B6").NumberFormat = $timeFmt
$ws10.Range("C6").Value = 0.77083333333333337
$ws10.Range("C6").NumberFormat = $timeFmt
$ws10.Range("D6").Value = "Created and tested view to display all orders of logged in user, added reduce, increase, and remove links on cart view, created view to display user info"
$ws10.Range("D6").WrapText = $true
$ws10.Range("E6").Value = 4
$ws10.Rows.Item(6).RowHeight = 39

# Row 7
$ws10.Range("A7").Value = 41712
$ws10.Range("A7").NumberFormat = $dateFmt
$ws10.Range("B7").Value = 0.95833333333333337
$ws10.Range("B7").NumberFormat = $timeFmt
$ws10.Range("C7").Value = 0.10416666666666667
$ws10.Range("C7").NumberFormat = $timeFmt
$ws10.Range("D7").Value = "Worked on restricting admin views from regular users"
$ws10.Range("D7").WrapText = $true
$ws10.Range("E7").Value = 3.5
$ws10.Rows.Item(7).RowHeight = 18

# ---------------------------------------------------------------------
# Week 11 - add rows 2 through 5
# ---------------------------------------------------------------------
$ws11 = $wb.Worksheets.Item("Week 11")

# Row 2
$ws11.Range("A2").Value = 41715
$ws11.Range("A2").NumberFormat = $dateFmt
$ws11.Range("B2").Value = 0.70833333333333337
$ws11.Range("B2").NumberFormat = $timeFmt
$ws11.Range("C2").Value = 0.72916666666666663
$ws11.Range("C2").NumberFormat = $timeFmt
$ws11.Range("D2").Value = "Added and tested user info update button"
$ws11.Range("E2").Value = 0.5
$ws11.Rows.Item(2).RowHeight = 18

# Row 3
$ws11.Range("A3").Value = 41716
$ws11.Range("A3").NumberFormat = $dateFmt
$ws11.Range("B3").Value = 0.91666666666666663
$ws11.Range("B3").NumberFormat = $timeFmt
$ws11.Range("C3").Value = 0.020833333333333332
$ws11.Range("C3").NumberFormat = $timeFmt
$ws11.Range("D3").Value = "Adjusted CSS on all pages"
$ws11.Range("E3").Value = 2.5
$ws11.Rows.Item(3).RowHeight = 18

# Row 4
$ws11.Range("A4").Value = 41719
$ws11.Range("A4").NumberFormat = $dateFmt
$ws11.Range("B4").Value = 0.625
$ws11.Range("B4").NumberFormat = $timeFmt
$ws11.Range("C4").Value = 0.75
$ws11.Range("C4").NumberFormat = $timeFmt
$ws11.Range("D4").Value = "Added incorrect password error message, created and tested password reset"
$ws11.Range("D4").WrapText = $true
$ws11.Range("E4").Value = 3
$ws11.Rows.Item(4).RowHeight = 26

# Row 5
$ws11.Range("A5").Value = 41720
$ws11.Range("A5").NumberFormat = $dateFmt
$ws11.Range("B5").Value = 0.5
$ws11.Range("B5").NumberFormat = $timeFmt
$ws11.Range("C5").Value = 0.58333333333333337
$ws11.Range("C5").NumberFormat = $timeFmt
$ws11.Range("D5").Value = "Worked on CSS and resposiveness"
$ws11.Range("E5").Value = 2
$ws11.Rows.Item(5).RowHeight = 18

# ---------------------------------------------------------------------
# Update selections and move the active tab from Week 10 to Week 11
# ---------------------------------------------------------------------
$ws10.Activate()
$ws10.Range("E8").Select()

$ws11.Activate()
$ws11.Range("A6").Select()
